# MagicCube.docx edit script
# Applies:
#  - proofErr spellcheck/grammar markers around a handful of runs
#    (splitting existing runs so the marked word sits in its own run)
#  - reflow of the "Hide the cursor / Keep focus / Drag red cross" tail of
#    the Small Task list: "Keep focus on the window." is replaced by
#    "Mouse can't move out of the window.", "Hide the cursor." and
#    "Drag red cross..." pick up strike-through, the _GoBack bookmark moves
#    up to sit on "Smaller the cube.", and "Drag red cross..." gets a
#    gramStart/gramEnd pair around "red cross".

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Title: "MagicCube." -> spellStart/gramStart ... spellEnd "." gramEnd
# ---------------------------------------------------------------------
$titleXml = @"
<w:p $wns>
  <w:pPr><w:pStyle w:val="Title"/></w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>MagicCube</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>.</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
$d.Paragraphs(1).Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) "Generate textureCoord." -> "Generate " + spellStart textureCoord spellEnd + "."
# ---------------------------------------------------------------------
$textureXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Generate </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>textureCoord</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(14).Range.InsertXML($textureXml)

# ---------------------------------------------------------------------
# 3) "Generate mipmap." -> "Generate " + spellStart mipmap spellEnd + "."
# ---------------------------------------------------------------------
$mipmapXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Generate </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>mipmap</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(19).Range.InsertXML($mipmapXml)

# ---------------------------------------------------------------------
# 4) "Camera can look&rotate to up/down." / "...to left/right."
#    -> "Camera can " + spellStart look&rotate spellEnd + " to .../...".
# ---------------------------------------------------------------------
$cameraUpDownXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Camera can </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>look&amp;rotate</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> to up/down.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(21).Range.InsertXML($cameraUpDownXml)

$cameraLeftRightXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Camera can </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>look&amp;rotate</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> to left/right.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(22).Range.InsertXML($cameraLeftRightXml)

# ---------------------------------------------------------------------
# 5) Reflow the tail of the Small Task list:
#    27 Smaller the cube.                     (keeps strike; gains bookmark)
#    28 Rotate the cube like a bullet.         (unchanged)
#    29 Find three new bullet textures. (...)  (unchanged)
#    30 Hide the cursor.                       (gains strike)
#    31 Keep focus on the window.  -->  Mouse can't move out of the window. (strike)
#    32 Drag red cross in the center...        (gains strike + gramStart/gramEnd)
# ---------------------------------------------------------------------
$tailXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>Smaller the cube.</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>Rotate the cube like a bullet</w:t></w:r>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>.</w:t></w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>Find three new bullet texture</w:t></w:r>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>s</w:t></w:r>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>. (Plastic, Wood, Metal)</w:t></w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>Hide the cursor.</w:t></w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>Mouse can&#8217;t move out of the window.</w:t></w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:strike/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Drag </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>red cross</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> in the center of the screen.</w:t></w:r>
</w:p>
"@

$startPara = $d.Paragraphs(27)
$endPara = $d.Paragraphs(32)
$tailRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$tailRange.InsertXML($tailXml)
